# Fill in the peer-assessment grades for the third team member (row 22,
# student 1220741) on the "Group and Self Assessment" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group and Self Assessment")

$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0

$excel.Calculate()
